$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 10385.333
$ws.Range("J18").Value = 12743
$ws.Range("L18").Value = 12743
$ws.Range("N18").Value = -13311
$ws.Range("H42").Value = 322.5
$ws.Range("J42").Value = 424.9091
$ws.Range("L42").Value = 1274.7273
$ws.Range("N42").Value = -1734.7273
$ws.Range("H51").Value = 8328
$ws.Range("J51").Value = 7497.5
$ws.Range("L51").Value = 7497.5
$ws.Range("N51").Value = -8465.5
$ws.Range("H112").Value = 3746.611
$ws.Range("J112").Value = 4763.3076
$ws.Range("L112").Value = 14289.9228
$ws.Range("N112").Value = -16505.9228
$ws.Range("H113").Value = 18024
$ws.Range("J113").Value = 10839
$ws.Range("L113").Value = 10839
$ws.Range("N113").Value = -17347
$ws.Range("H125").Value = 4271.5454
$ws.Range("I125").Value = 5174.4
$ws.Range("K125").Value = 46569.6
$ws.Range("M125").Value = -44109.6
$ws.Range("H133").Value = 80809.60000000001
$ws.Range("J133").Value = 80809.60000000001
$ws.Range("L133").Value = 80809.60000000001
$ws.Range("N133").Value = -90929.60000000001
$ws.Range("H137").Value = 19482.416
$ws.Range("I137").Value = 27024.25
$ws.Range("J137").Value = 4398.75
$ws.Range("K137").Value = 81072.75
$ws.Range("L137").Value = 13196.25
$ws.Range("M137").Value = -78522.75
$ws.Range("N137").Value = -18296.25
$ws.Range("H138").Value = 3267.516
$ws.Range("J138").Value = 4384.5
$ws.Range("L138").Value = 13153.5
$ws.Range("N138").Value = -23433.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 10003
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 10003
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10003
$ws.Range("M6").ClearContents()
$ws.Range("N6").Value = -10349
$ws.Range("H10").Value = 1263
$ws.Range("I10").Value = 1263
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 1263
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -1093
$ws.Range("N10").ClearContents()
$ws.Range("H12").Value = 151.5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 151.5
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 151.5
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -497.5
$ws.Range("H32").Value = 6147.04
$ws.Range("I32").Value = 5986.5
$ws.Range("K32").Value = 5986.5
$ws.Range("M32").Value = -5699.5
$ws.Range("H61").Value = 5349.6
$ws.Range("I61").Value = 5360.0415
$ws.Range("J61").Value = 5099
$ws.Range("K61").Value = 5360.0415
$ws.Range("L61").Value = 5099
$ws.Range("M61").Value = -5148.0415
$ws.Range("N61").Value = -5523
$ws.Range("H74").Value = 5964.8667
$ws.Range("I74").Value = 1947.3
$ws.Range("K74").Value = 1947.3
$ws.Range("M74").Value = -1073.3
$ws.Range("H77").Value = 5964.8667
$ws.Range("I77").Value = 1947.3
$ws.Range("K77").Value = 9736.5
$ws.Range("M77").Value = -5368.5
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820
$ws.Range("H132").Value = 5031.321
$ws.Range("I132").Value = 5494.657
$ws.Range("J132").Value = 4130.3887
$ws.Range("K132").Value = 16483.971
$ws.Range("L132").Value = 12391.1661
$ws.Range("M132").Value = -13953.971
$ws.Range("N132").Value = -17451.1661
$ws.Range("H136").Value = 5349.6
$ws.Range("I136").Value = 5360.0415
$ws.Range("J136").Value = 5099
$ws.Range("K136").Value = 16080.1245
$ws.Range("L136").Value = 15297
$ws.Range("M136").Value = -13530.1245
$ws.Range("N136").Value = -20397
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H99").Value = 18767.459
$ws.Range("I99").Value = 29884.076
$ws.Range("J99").Value = 5629.636
$ws.Range("K99").Value = 29884.076
$ws.Range("L99").Value = 5629.636
$ws.Range("M99").Value = -28386.076
$ws.Range("N99").Value = -8625.636
$ws.Range("H134").Value = 14124.952
$ws.Range("I134").Value = 15032.842
$ws.Range("J134").Value = 5500
$ws.Range("K134").Value = 45098.526
$ws.Range("L134").Value = 16500
$ws.Range("M134").Value = -42563.526
$ws.Range("N134").Value = -21570
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 509
$ws.Range("I5").Value = 331.42856
$ws.Range("J5").Value = 923.3333
$ws.Range("K5").Value = 331.42856
$ws.Range("L5").Value = 923.3333
$ws.Range("M5").Value = -219.42856
$ws.Range("N5").Value = -1147.3333
$ws.Range("H6").Value = 2754664.8
$ws.Range("I6").Value = 4125499.8
$ws.Range("K6").Value = 4125499.8
$ws.Range("M6").Value = -4125386.8
$ws.Range("H7").Value = 8920.434999999999
$ws.Range("I7").Value = 12643.75
$ws.Range("K7").Value = 12643.75
$ws.Range("M7").Value = -12530.75
$ws.Range("H31").Value = 7623.56
$ws.Range("I31").Value = 9724.3125
$ws.Range("K31").Value = 9724.3125
$ws.Range("M31").Value = -9429.3125
$ws.Range("H34").Value = 7623.56
$ws.Range("I34").Value = 9724.3125
$ws.Range("K34").Value = 9724.3125
$ws.Range("M34").Value = -9522.3125
$ws.Range("H58").Value = 3947.8
$ws.Range("I58").Value = 5179.778
$ws.Range("K58").Value = 5179.778
$ws.Range("M58").Value = -4976.778
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -104910
$ws.Range("H134").Value = 3661.1428
$ws.Range("I134").Value = 4105.9
$ws.Range("K134").Value = 12317.7
$ws.Range("M134").Value = -9782.699999999999
$ws.Range("H136").Value = 3947.8
$ws.Range("I136").Value = 5179.778
$ws.Range("K136").Value = 15539.334
$ws.Range("M136").Value = -12989.334
$ws.Range("H141").Value = 255939.19
$ws.Range("J141").Value = 301879.62
$ws.Range("L141").Value = 301879.62
$ws.Range("N141").Value = -312239.62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 417491.72
$ws.Range("I5").Value = 275.14285
$ws.Range("J5").Value = 1001594.9
$ws.Range("K5").Value = 825.4285500000001
$ws.Range("L5").Value = 3004784.7
$ws.Range("M5").Value = -713.4285500000001
$ws.Range("N5").Value = -3005008.7
$ws.Range("H122").Value = 4861.8857
$ws.Range("J122").Value = 5421.033
$ws.Range("L122").Value = 48789.29700000001
$ws.Range("N122").Value = -53689.29700000001
$ws.Range("H135").Value = 417491.72
$ws.Range("I135").Value = 275.14285
$ws.Range("J135").Value = 1001594.9
$ws.Range("K135").Value = 2476.28565
$ws.Range("L135").Value = 9014354.1
$ws.Range("M135").Value = 58.71434999999974
$ws.Range("N135").Value = -9019424.1
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8198.823
$ws.Range("I70").Value = 6467.0835
$ws.Range("K70").Value = 6467.0835
$ws.Range("M70").Value = -6197.0835
$ws.Range("H73").Value = 8198.823
$ws.Range("I73").Value = 6467.0835
$ws.Range("K73").Value = 6467.0835
$ws.Range("M73").Value = -5531.0835
$ws.Range("H94").Value = 40853708
$ws.Range("I94").Value = 1000000
$ws.Range("J94").Value = 44839080
$ws.Range("K94").Value = 1000000
$ws.Range("L94").Value = 44839080
$ws.Range("M94").Value = -999324
$ws.Range("N94").Value = -44840432
$ws.Range("H122").Value = 6627.2896
$ws.Range("I122").Value = 4017.516
$ws.Range("K122").Value = 12052.548
$ws.Range("M122").Value = -9602.548000000001
$ws.Range("H123").Value = 49777
$ws.Range("J123").Value = 49777
$ws.Range("L123").Value = 49777
$ws.Range("N123").Value = -54677
$ws.Range("H126").Value = 8386.606
$ws.Range("I126").Value = 9837.8125
$ws.Range("J126").Value = 7020.7646
$ws.Range("K126").Value = 29513.4375
$ws.Range("L126").Value = 21062.2938
$ws.Range("M126").Value = -27043.4375
$ws.Range("N126").Value = -26002.2938
$ws.Range("H132").Value = 2593.0667
$ws.Range("I132").Value = 2375.5
$ws.Range("J132").Value = 3463.3333
$ws.Range("K132").Value = 7126.5
$ws.Range("L132").Value = 10389.9999
$ws.Range("M132").Value = -4596.5
$ws.Range("N132").Value = -15449.9999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19870.066
$ws.Range("I7").Value = 35408.92
$ws.Range("K7").Value = 35408.92
$ws.Range("M7").Value = -35296.92
$ws.Range("H126").Value = 19870.066
$ws.Range("I126").Value = 35408.92
$ws.Range("K126").Value = 106226.76
$ws.Range("M126").Value = -103756.76
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 10625.29
$ws.Range("I107").Value = 1011.8461
$ws.Range("J107").Value = 60615.2
$ws.Range("K107").Value = 3035.5383
$ws.Range("L107").Value = 181845.6
$ws.Range("M107").Value = -1115.5383
$ws.Range("N107").Value = -185685.6
$ws.Range("H122").Value = 4179.4653
$ws.Range("I122").Value = 1565.4482
$ws.Range("J122").Value = 9594.214
$ws.Range("K122").Value = 4696.3446
$ws.Range("L122").Value = 28782.642
$ws.Range("M122").Value = -2246.3446
$ws.Range("N122").Value = -33682.642
$ws.Range("H126").Value = 30292.666
$ws.Range("I126").Value = 42209.7
$ws.Range("J126").Value = 6458.6
$ws.Range("K126").Value = 126629.1
$ws.Range("L126").Value = 19375.8
$ws.Range("M126").Value = -124159.1
$ws.Range("N126").Value = -24315.8
$ws.Range("H136").Value = 360626.12
$ws.Range("I136").Value = 552020.9399999999
$ws.Range("J136").Value = 3355.8667
$ws.Range("K136").Value = 1656062.82
$ws.Range("L136").Value = 10067.6001
$ws.Range("M136").Value = -1653512.82
$ws.Range("N136").Value = -15167.6001
$ws.Range("H137").Value = 104921.664
$ws.Range("J137").Value = 104921.664
$ws.Range("L137").Value = 104921.664
$ws.Range("N137").Value = -115121.664
